$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 94, shifting existing rows 94..190 down to 95..191
$ws.Rows.Item(94).Insert([Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftDown)

# Fill in the new row 94 with values.
# Columns A, B, C, E, F, G, H, I, N, O, Q, R are copied from the row that
# used to be row 94 (now row 95, identical to the original row94 before the shift).
$ws.Cells.Item(94, 1).Value = $ws.Cells.Item(95, 1).Value()
$ws.Cells.Item(94, 2).Value = $ws.Cells.Item(95, 2).Value()
$ws.Cells.Item(94, 3).Value = $ws.Cells.Item(95, 3).Value()
$ws.Cells.Item(94, 4).Value = 44629
$ws.Cells.Item(94, 4).NumberFormat = $ws.Cells.Item(95, 4).NumberFormat
$ws.Cells.Item(94, 5).Value = $ws.Cells.Item(95, 5).Value()
$ws.Cells.Item(94, 6).Value = $ws.Cells.Item(95, 6).Value()
$ws.Cells.Item(94, 7).Value = $ws.Cells.Item(95, 7).Value()
$ws.Cells.Item(94, 8).Value = $ws.Cells.Item(95, 8).Value()
$ws.Cells.Item(94, 9).Value = $ws.Cells.Item(95, 9).Value()
$ws.Cells.Item(94, 10).Value = 300
$ws.Cells.Item(94, 11).Value = 1100
$ws.Cells.Item(94, 12).Value = 1200
$ws.Cells.Item(94, 13).Value = 1150
$ws.Cells.Item(94, 14).Value = $ws.Cells.Item(95, 14).Value()
$ws.Cells.Item(94, 15).Value = $ws.Cells.Item(95, 15).Value()
$ws.Cells.Item(94, 16).Value = 1150
$ws.Cells.Item(94, 17).Value = $ws.Cells.Item(95, 17).Value()
$ws.Cells.Item(94, 18).Value = $ws.Cells.Item(95, 18).Value()
